# Commit: "merged CDS test suites to create CDS_Regression suite"
#
# This applies two content edits to the "startup" sheet of the workbook:
#   1) Cell B2 holds a Cypher query (shared string). Update it so that:
#        - the sample_id collection no longer falls back to the literal
#          "Not specified in data" text (apoc.coll.sort(collect(distinct
#          coalesce(samp.sample_id, "Not specified in data"))) ->
#          apoc.coll.sort(collect(distinct samp.sample_id)))
#        - the trailing "ORDER BY p.participant_id LIMIT 100" is
#          lower-cased to "ORDER BY p.participant_id limit 100"
#   2) The active selection on the sheet moves from D3 to D2.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

# --- 1) Update the Cypher query text stored in B2 -------------------------
$query = $ws.Range("B2").Value2

$oldFragment = 'apoc.coll.sort(collect(distinct coalesce(samp.sample_id, "Not specified in data"))) as samp'
$newFragment = 'apoc.coll.sort(collect(distinct samp.sample_id)) as samp'
$query = $query.Replace($oldFragment, $newFragment)

$oldTail = 'ORDER BY p.participant_id LIMIT 100'
$newTail = 'ORDER BY p.participant_id limit 100'
$query = $query.Replace($oldTail, $newTail)

$ws.Range("B2").Value = $query

# --- 2) Move the active cell / selection from D3 to D2 ---------------------
$ws.Range("D2").Select()
